$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row right below the header (new row 2), shifting all
# existing data rows down by one (old row 2 -> 3, ..., old row 24 -> 25).
$ws.Rows.Item(2).Insert()

# The freshly inserted row inherits formatting from the row above (the bold,
# bordered header style). Reset it back to the plain "Normal" style used by
# every other data row before re-applying the date format to column D.
$ws.Rows.Item(2).Style = "Normal"
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Populate the new row with the latest weekly price record.
$ws.Range("A2").Value = 10
$ws.Range("B2").Value = "Vega Modelo de Temuco"
$ws.Range("C2").Value = "La Araucanía"
$ws.Range("D2").Value = 44956
$ws.Range("E2").Value = 9
$ws.Range("F2").Value = 100112017
$ws.Range("G2").Value = "Ramas de apio"
$ws.Range("H2").Value = "Sin especificar"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 40
$ws.Range("K2").Value = 5000
$ws.Range("L2").Value = 5000
$ws.Range("M2").Value = 5000
$ws.Range("N2").Value = "`$/paquete"
$ws.Range("O2").Value = "Región de La Araucanía"
$ws.Range("P2").Value = 5000
$ws.Range("Q2").Value = 1
$ws.Range("R2").Value = "Hortaliza"
